$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value is a plain decimal number would be
# auto-converted from text to a numeric type by Excel's smart-entry logic.
# Force those specific cells to Text format first so the literal digit
# string (including trailing zeros) is preserved exactly, matching the
# original inline-string cell content.
$textCells = @(
    "D4",
    "D5",
    "D6",
    "D7",
    "D9",
    "D10",
    "D11",
    "D12",
    "D13",
    "D14",
    "D19",
    "D20",
    "D21",
    "D22",
    "D23",
    "D24",
    "D25",
    "D26",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D41",
    "D42",
    "D43",
    "D44",
    "D46",
    "D47",
    "D48",
    "D49",
    "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "63.124.62"
$ws.Range("E2").Value = "  -1.42%  "

# Row 3
$ws.Range("D3").Value = "3.083.74"
$ws.Range("E3").Value = "  +0.49%  "

# Row 4
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.19%  "

# Row 5
$ws.Range("D5").Value = "555.60"
$ws.Range("E5").Value = "  +0.87%  "

# Row 6
$ws.Range("D6").Value = "137.66"
$ws.Range("E6").Value = "  -2.34%  "

# Row 7
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.17%  "

# Row 8
$ws.Range("D8").Value = "3.080.66"
$ws.Range("E8").Value = "  +0.69%  "

# Row 9
$ws.Range("D9").Value = "0.495"
$ws.Range("E9").Value = "  +1.20%  "

# Row 10
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "0.163"
$ws.Range("E10").Value = "  +7.05%  "

# Row 11
$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D11").Value = "6.70"
$ws.Range("E11").Value = "  +2.28%  "

# Row 12
$ws.Range("D12").Value = "0.453"
$ws.Range("E12").Value = "  +0.08%  "

# Row 13
$ws.Range("D13").Value = "35.27"
$ws.Range("E13").Value = "  -1.29%  "

# Row 14
$ws.Range("D14").Value = "0.0000218"
$ws.Range("E14").Value = "  +0.58%  "

# Row 15
$ws.Range("D15").Value = "3.579.91"
$ws.Range("E15").Value = "  +0.52%  "

# Row 16
$ws.Range("D16").Value = "63.044.99"
$ws.Range("E16").Value = "  -1.59%  "

# Row 17
$ws.Range("E17").Value = "  -0.17%  "

# Row 18
$ws.Range("D18").Value = "3.068.95"
$ws.Range("E18").Value = "  -0.20%  "

# Row 19
$ws.Range("D19").Value = "506.95"
$ws.Range("E19").Value = "  +4.11%  "

# Row 20
$ws.Range("D20").Value = "6.65"
$ws.Range("E20").Value = "  +0.72%  "

# Row 21
$ws.Range("D21").Value = "13.63"
$ws.Range("E21").Value = "  +0.43%  "

# Row 22
$ws.Range("D22").Value = "0.706"
$ws.Range("E22").Value = "  +3.18%  "

# Row 23
$ws.Range("D23").Value = "7.27"
$ws.Range("E23").Value = "  +0.55%  "

# Row 24
$ws.Range("D24").Value = "77.68"
$ws.Range("E24").Value = "  -0.35%  "

# Row 25
$ws.Range("D25").Value = "12.34"
$ws.Range("E25").Value = "  -0.35%  "

# Row 26
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.35%  "

# Row 27
$ws.Range("D27").Value = "2.77"
$ws.Range("E27").Value = "  +2.13%  "

# Row 28
$ws.Range("D28").Value = "8.33"
$ws.Range("E28").Value = "  +1.39%  "

# Row 29
$ws.Range("D29").Value = "2.04"
$ws.Range("E29").Value = "  -1.60%  "

# Row 30
$ws.Range("D30").Value = "0.996"
$ws.Range("E30").Value = "  -0.48%  "

# Row 31
$ws.Range("D31").Value = "26.23"
$ws.Range("E31").Value = "  +2.01%  "

# Row 32
$ws.Range("D32").Value = "2.52"
$ws.Range("E32").Value = "  -4.71%  "

# Row 33
$ws.Range("D33").Value = "1.11"
$ws.Range("E33").Value = "  -3.01%  "

# Row 34
$ws.Range("D34").Value = "535.02"
$ws.Range("E34").Value = "  -9.73%  "

# Row 35
$ws.Range("D35").Value = "57.68"
$ws.Range("E35").Value = "  +10.91%  "

# Row 36
$ws.Range("D36").Value = "5.90"
$ws.Range("E36").Value = "  -0.83%  "

# Row 37
$ws.Range("D37").Value = "5.16"
$ws.Range("E37").Value = "  -4.38%  "

# Row 38
$ws.Range("D38").Value = "0.0413"
$ws.Range("E38").Value = "  +3.80%  "

# Row 39
$ws.Range("D39").Value = "0.0796"
$ws.Range("E39").Value = "  +0.85%  "

# Row 40
$ws.Range("D40").Value = "3.065.82"
$ws.Range("E40").Value = "  +3.33%  "

# Row 41
$ws.Range("D41").Value = "0.118"
$ws.Range("E41").Value = "  +0.09%  "

# Row 42
$ws.Range("D42").Value = "8.10"
$ws.Range("E42").Value = "  -1.05%  "

# Row 43
$ws.Range("D43").Value = "2.66"
$ws.Range("E43").Value = "  -7.18%  "

# Row 44
$ws.Range("D44").Value = "0.254"
$ws.Range("E44").Value = "  +4.10%  "

# Row 45
$ws.Range("E45").Value = "  +0.03%  "

# Row 46
$ws.Range("D46").Value = "2.08"
$ws.Range("E46").Value = "  -0.17%  "

# Row 47
$ws.Range("D47").Value = "120.91"
$ws.Range("E47").Value = "  +0.91%  "

# Row 48
$ws.Range("D48").Value = "24.18"
$ws.Range("E48").Value = "  -2.92%  "

# Row 49
$ws.Range("D49").Value = "0.107"
$ws.Range("E49").Value = "  -0.31%  "

# Row 50
$ws.Range("D50").Value = "0.0₃0495"
$ws.Range("E50").Value = "  -6.23%  "

# Row 51
$ws.Range("B51").Value = "CoreDAO"
$ws.Range("C51").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D51").Value = "2.35"
$ws.Range("E51").Value = "  +66.66%  "
